$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after the dbRDA* row (row 10), pushing the legend
# block (rows 12-20) down to rows 14-22.
$ws.Rows("11:12").Insert()

# The row-insert copies the formatting (incl. number formats) of row 10
# across columns B:P into the two new rows; clear that so the new rows
# only keep the column-A "Good" style, matching the target layout.
$ws.Range("B11:P12").Clear()

# Populate the two new method rows (added in shared-string order so the
# resulting shared string table indices line up with the target file).
$ws.Range("A11").Value = "PERMDISP"
$ws.Range("A12").Value = "BEST"

# Rename the dbRDA* method to db-RDA* (added to the shared string table
# last, after the two new rows above).
$ws.Range("A10").Value = "db-RDA*"

# Match the author's final selection in the worksheet.
[void]$ws.Range("B12").Select()
